$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 997.58826
$ws.Range("J32").Value = 1142.3
$ws.Range("L32").Value = 1142.3
$ws.Range("N32").Value = -1794.3

$ws.Range("H33").Value = 682.6667
$ws.Range("I33").Value = 519.2
$ws.Range("K33").Value = 519.2
$ws.Range("M33").Value = -290.2

$ws.Range("H70").Value = 4384.615
$ws.Range("I70").Value = 3500
$ws.Range("K70").Value = 10500
$ws.Range("M70").Value = -10230

$ws.Range("H73").Value = 4384.615
$ws.Range("I73").Value = 3500
$ws.Range("K73").Value = 10500
$ws.Range("M73").Value = -9564

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 52000
$ws.Range("J29").Value = 52000
$ws.Range("L29").Value = 52000
$ws.Range("N29").Value = -52616

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3418.8333
$ws.Range("I20").Value = 3800.8
$ws.Range("K20").Value = 3800.8
$ws.Range("M20").Value = -3553.8

$ws.Range("H22").Value = 10146.25
$ws.Range("I22").Value = 197
$ws.Range("J22").Value = 20095.5
$ws.Range("K22").Value = 197
$ws.Range("L22").Value = 20095.5
$ws.Range("M22").Value = -24
$ws.Range("N22").Value = -20441.5

$ws.Range("H54").Value = 29999.5
$ws.Range("I54").Value = 29999.5
$ws.Range("K54").Value = 29999.5
$ws.Range("M54").Value = -29515.5

$ws.Range("H99").Value = 4581.8
$ws.Range("I99").Value = 4652.5
$ws.Range("K99").Value = 4652.5
$ws.Range("M99").Value = -3154.5

$ws.Range("H105").Value = 3210.1875
$ws.Range("I105").Value = 3087.9092
$ws.Range("J105").Value = 3479.2
$ws.Range("K105").Value = 3087.9092
$ws.Range("L105").Value = 3479.2
$ws.Range("M105").Value = -1340.9092
$ws.Range("N105").Value = -6973.2

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 2163.3333
$ws.Range("I32").Value = 999.5
$ws.Range("K32").Value = 999.5
$ws.Range("M32").Value = -683.5

$ws.Range("H58").Value = 3227.4285
$ws.Range("I58").Value = 2118.4
$ws.Range("K58").Value = 2118.4
$ws.Range("M58").Value = -1915.4

$ws.Range("H99").Value = 4263.4
$ws.Range("J99").Value = 4940.778
$ws.Range("L99").Value = 4940.778
$ws.Range("N99").Value = -7936.778

$ws.Range("H105").Value = 1305.5714
$ws.Range("I105").Value = 1305.5714
$ws.Range("K105").Value = 1305.5714
$ws.Range("M105").Value = 441.4286

$ws.Range("H126").Value = 4263.4
$ws.Range("J126").Value = 4940.778
$ws.Range("L126").Value = 14822.334
$ws.Range("N126").Value = -19762.334

$ws.Range("H132").Value = 2346.524
$ws.Range("I132").Value = 1517.3125
$ws.Range("K132").Value = 4551.9375
$ws.Range("M132").Value = -2021.9375

$ws.Range("H134").Value = 1246
$ws.Range("I134").Value = 1038.25
$ws.Range("J134").Value = 1800
$ws.Range("K134").Value = 3114.75
$ws.Range("L134").Value = 5400
$ws.Range("M134").Value = -579.75
$ws.Range("N134").Value = -10470

$ws.Range("H136").Value = 3227.4285
$ws.Range("I136").Value = 2118.4
$ws.Range("K136").Value = 6355.200000000001
$ws.Range("M136").Value = -3805.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H134").Value = 171821.83
$ws.Range("I134").Value = 333632.66
$ws.Range("J134").Value = 10011
$ws.Range("K134").Value = 1000897.98
$ws.Range("L134").Value = 30033
$ws.Range("M134").Value = -995827.98
$ws.Range("N134").Value = -40173

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 289.16666
$ws.Range("I2").Value = 303.8
$ws.Range("K2").Value = 303.8
$ws.Range("M2").Value = -190.8

$ws.Range("H3").Value = 40399.375
$ws.Range("J3").Value = 52900
$ws.Range("L3").Value = 52900
$ws.Range("N3").Value = -53132

$ws.Range("H21").Value = 26000
$ws.Range("J21").Value = 26000
$ws.Range("L21").Value = 26000
$ws.Range("N21").Value = -26346

$ws.Range("H24").Value = 4413888.5

$ws.Range("H30").Value = 26000
$ws.Range("J30").Value = 26000
$ws.Range("L30").Value = 26000
$ws.Range("N30").Value = -26210

$ws.Range("H39").Value = 30000
$ws.Range("J39").Value = 30000
$ws.Range("L39").Value = 30000
$ws.Range("N39").Value = -31064

$ws.Range("H70").Value = 333333340
$ws.Range("I70").Value = 333333340
$ws.Range("K70").Value = 333333340
$ws.Range("M70").Value = -333333070

$ws.Range("H73").Value = 333333340
$ws.Range("I73").Value = 333333340
$ws.Range("K73").Value = 333333340
$ws.Range("M73").Value = -333332404

$ws.Range("H100").Value = 38900
$ws.Range("J100").Value = 38900
$ws.Range("L100").Value = 38900
$ws.Range("N100").Value = -41064

$ws.Range("H113").Value = 3271.3333
$ws.Range("I113").Value = 2512.5
$ws.Range("K113").Value = 2512.5
$ws.Range("M113").Value = -342.5

$ws.Range("H126").Value = 16633.166
$ws.Range("I126").Value = 16633.166
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 49899.49800000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -47429.49800000001
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19713.188
$ws.Range("I7").Value = 19210.732
$ws.Range("K7").Value = 19210.732
$ws.Range("M7").Value = -19098.732

$ws.Range("H22").Value = 1284.4286
$ws.Range("J22").Value = 1497.75
$ws.Range("L22").Value = 1497.75
$ws.Range("N22").Value = -2087.75

$ws.Range("H23").Value = 2499.5
$ws.Range("J23").Value = 2499.5
$ws.Range("L23").Value = 2499.5
$ws.Range("N23").Value = -2959.5

$ws.Range("H27").Value = 1284.4286
$ws.Range("J27").Value = 1497.75
$ws.Range("L27").Value = 1497.75
$ws.Range("N27").Value = -1711.75

$ws.Range("H55").Value = 316
$ws.Range("J55").Value = 329.25
$ws.Range("L55").Value = 329.25
$ws.Range("N55").Value = -675.25

$ws.Range("H123").Value = 77997
$ws.Range("J123").Value = 77997
$ws.Range("L123").Value = 77997
$ws.Range("N123").Value = -87797

$ws.Range("H126").Value = 19713.188
$ws.Range("I126").Value = 19210.732
$ws.Range("K126").Value = 57632.196
$ws.Range("M126").Value = -55162.196

$ws.Range("H132").Value = 5062.4
$ws.Range("I132").Value = 3981
$ws.Range("J132").Value = 9388
$ws.Range("K132").Value = 11943
$ws.Range("L132").Value = 28164
$ws.Range("M132").Value = -9413
$ws.Range("N132").Value = -33224

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 10029
$ws.Range("J32").Value = 10029
$ws.Range("L32").Value = 10029
$ws.Range("N32").Value = -10663

$ws.Range("H34").Value = 32500
$ws.Range("I34").Value = 40000
$ws.Range("J34").Value = 25000
$ws.Range("K34").Value = 40000
$ws.Range("L34").Value = 25000
$ws.Range("M34").Value = -39797
$ws.Range("N34").Value = -25406

$ws.Range("H107").Value = 2244.25
$ws.Range("I107").Value = 2084.75
$ws.Range("K107").Value = 6254.25
$ws.Range("M107").Value = -4334.25

$ws.Range("H132").Value = 3316.7273
$ws.Range("I132").Value = 3165
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 9495
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -6965
$ws.Range("N132").Value = -17058.5

$ws.Range("H136").Value = 3612.3333
$ws.Range("I136").Value = 3612.3333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10836.9999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8286.999899999999
$ws.Range("N136").ClearContents()

$ws.Range("H139").Value = 60250
$ws.Range("J139").Value = 60250
$ws.Range("L139").Value = 60250
$ws.Range("N139").Value = -70530
